$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 898.7143
$ws.Range("I40").Value = 818
$ws.Range("J40").Value = 1100.5
$ws.Range("K40").Value = 818
$ws.Range("L40").Value = 1100.5
$ws.Range("M40").Value = -643
$ws.Range("N40").Value = -1450.5

$ws.Range("H132").Value = 1363.3877
$ws.Range("I132").Value = 1144.881
$ws.Range("J132").Value = 2674.4285
$ws.Range("K132").Value = 3434.643
$ws.Range("L132").Value = 8023.2855
$ws.Range("M132").Value = -904.643
$ws.Range("N132").Value = -13083.2855

$ws.Range("H135").Value = 3477.7576
$ws.Range("I135").Value = 2737.2
$ws.Range("J135").Value = 5792
$ws.Range("K135").Value = 24634.8
$ws.Range("L135").Value = 52128
$ws.Range("M135").Value = -22099.8
$ws.Range("N135").Value = -57198

$ws.Range("H136").Value = 43000
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 43000
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 43000
$ws.Range("N136").Value = -53200

$ws.Range("H137").Value = 1714.1666
$ws.Range("I137").Value = 1135.6
$ws.Range("J137").Value = 2437.375
$ws.Range("K137").Value = 3406.8
$ws.Range("L137").Value = 7312.125
$ws.Range("M137").Value = -856.7999999999997
$ws.Range("N137").Value = -12412.125

$ws.Range("H138").Value = 4831.593
$ws.Range("I138").Value = 3781.0588
$ws.Range("J138").Value = 5314.2705
$ws.Range("K138").Value = 11343.1764
$ws.Range("L138").Value = 15942.8115
$ws.Range("M138").Value = -6203.1764
$ws.Range("N138").Value = -26222.8115

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 616565.5600000001
$ws.Range("I32").Value = 9689.24
$ws.Range("J32").Value = 1995830
$ws.Range("K32").Value = 9689.24
$ws.Range("L32").Value = 1995830
$ws.Range("M32").Value = -9402.24
$ws.Range("N32").Value = -1996404

$ws.Range("H37").Value = 6424.7
$ws.Range("I37").Value = 2344.6667
$ws.Range("J37").Value = 8173.2856
$ws.Range("K37").Value = 2344.6667
$ws.Range("L37").Value = 8173.2856
$ws.Range("M37").Value = -2071.6667
$ws.Range("N37").Value = -8719.285599999999

$ws.Range("H61").Value = 3675.0667
$ws.Range("I61").Value = 3643.7144
$ws.Range("J61").Value = 4114
$ws.Range("K61").Value = 3643.7144
$ws.Range("L61").Value = 4114
$ws.Range("M61").Value = -3431.7144
$ws.Range("N61").Value = -4538

$ws.Range("H74").Value = 1053.0769
$ws.Range("I74").Value = 1056.3334
$ws.Range("J74").Value = 1014
$ws.Range("K74").Value = 1056.3334
$ws.Range("L74").Value = 1014
$ws.Range("M74").Value = -182.3334
$ws.Range("N74").Value = -2762

$ws.Range("H77").Value = 1053.0769
$ws.Range("I77").Value = 1056.3334
$ws.Range("J77").Value = 1014
$ws.Range("K77").Value = 5281.666999999999
$ws.Range("L77").Value = 5070
$ws.Range("M77").Value = -913.6669999999995
$ws.Range("N77").Value = -13806

$ws.Range("H136").Value = 3675.0667
$ws.Range("I136").Value = 3643.7144
$ws.Range("J136").Value = 4114
$ws.Range("K136").Value = 10931.1432
$ws.Range("L136").Value = 12342
$ws.Range("M136").Value = -8381.143199999999
$ws.Range("N136").Value = -17442

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 6215.136
$ws.Range("I82").Value = 2094.8462
$ws.Range("J82").Value = 12166.667
$ws.Range("K82").Value = 2094.8462
$ws.Range("L82").Value = 12166.667
$ws.Range("M82").Value = -1711.8462
$ws.Range("N82").Value = -12932.667

$ws.Range("H85").Value = 6215.136
$ws.Range("I85").Value = 2094.8462
$ws.Range("J85").Value = 12166.667
$ws.Range("K85").Value = 2094.8462
$ws.Range("L85").Value = 12166.667
$ws.Range("M85").Value = -768.8462
$ws.Range("N85").Value = -14818.667

$ws.Range("H107").Value = 5105
$ws.Range("I107").Value = 7230.5
$ws.Range("J107").Value = 1704.2
$ws.Range("K107").Value = 7230.5
$ws.Range("L107").Value = 1704.2
$ws.Range("M107").Value = -5310.5
$ws.Range("N107").Value = -5544.2

$ws.Range("H134").Value = 18853
$ws.Range("I134").Value = 927.75
$ws.Range("J134").Value = 36778.25
$ws.Range("K134").Value = 2783.25
$ws.Range("L134").Value = 110334.75
$ws.Range("M134").Value = -248.25
$ws.Range("N134").Value = -115404.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10298.632
$ws.Range("I31").Value = 3468.6128
$ws.Range("J31").Value = 15003.756
$ws.Range("K31").Value = 3468.6128
$ws.Range("L31").Value = 15003.756
$ws.Range("M31").Value = -3173.6128
$ws.Range("N31").Value = -15593.756

$ws.Range("H34").Value = 10298.632
$ws.Range("I34").Value = 3468.6128
$ws.Range("J34").Value = 15003.756
$ws.Range("K34").Value = 3468.6128
$ws.Range("L34").Value = 15003.756
$ws.Range("M34").Value = -3266.6128
$ws.Range("N34").Value = -15407.756

$ws.Range("H50").Value = 8639.286
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 8639.286
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 8639.286
$ws.Range("N50").Value = -9889.286

$ws.Range("H58").Value = 1401.5
$ws.Range("I58").Value = 839.2143
$ws.Range("J58").Value = 2188.7
$ws.Range("K58").Value = 839.2143
$ws.Range("L58").Value = 2188.7
$ws.Range("M58").Value = -636.2143
$ws.Range("N58").Value = -2594.7

$ws.Range("H59").Value = 13842
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 13842
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 13842
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -16132

$ws.Range("H60").Value = 6899.625
$ws.Range("I60").Value = 5498.25
$ws.Range("J60").Value = 8301
$ws.Range("K60").Value = 5498.25
$ws.Range("L60").Value = 8301
$ws.Range("M60").Value = -4987.25
$ws.Range("N60").Value = -9323

$ws.Range("H68").Value = 16220.875
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 16220.875
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 16220.875
$ws.Range("N68").Value = -17718.875

$ws.Range("H71").Value = 16220.875
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 16220.875
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 48662.625
$ws.Range("N71").Value = -56150.625

$ws.Range("H74").Value = 17476.857
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 17476.857
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 17476.857
$ws.Range("N74").Value = -19224.857

$ws.Range("H77").Value = 17476.857
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 17476.857
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 52430.571
$ws.Range("N77").Value = -61166.571

$ws.Range("H86").Value = 30546.482
$ws.Range("I86").Value = 5779.2354
$ws.Range("J86").Value = 72650.8
$ws.Range("K86").Value = 5779.2354
$ws.Range("L86").Value = 72650.8
$ws.Range("M86").Value = -4656.2354
$ws.Range("N86").Value = -74896.8

$ws.Range("H89").Value = 30546.482
$ws.Range("I89").Value = 5779.2354
$ws.Range("J89").Value = 72650.8
$ws.Range("K89").Value = 28896.177
$ws.Range("L89").Value = 363254
$ws.Range("M89").Value = -23280.177
$ws.Range("N89").Value = -374486

$ws.Range("H136").Value = 1401.5
$ws.Range("I136").Value = 839.2143
$ws.Range("J136").Value = 2188.7
$ws.Range("K136").Value = 2517.6429
$ws.Range("L136").Value = 6566.099999999999
$ws.Range("M136").Value = 32.35710000000017
$ws.Range("N136").Value = -11666.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 1481.1111
$ws.Range("I80").Value = 500
$ws.Range("J80").Value = 1538.8235
$ws.Range("K80").Value = 1500
$ws.Range("L80").Value = 4616.470499999999
$ws.Range("M80").Value = -564
$ws.Range("N80").Value = -6488.470499999999

$ws.Range("H83").Value = 1481.1111
$ws.Range("I83").Value = 500
$ws.Range("J83").Value = 1538.8235
$ws.Range("K83").Value = 4500
$ws.Range("L83").Value = 13849.4115
$ws.Range("M83").Value = 180
$ws.Range("N83").Value = -23209.4115

$ws.Range("H117").Value = 92781.82000000001
$ws.Range("I117").Value = 1350
$ws.Range("J117").Value = 113100
$ws.Range("K117").Value = 4050
$ws.Range("L117").Value = 339300
$ws.Range("M117").Value = -608
$ws.Range("N117").Value = -346184

$ws.Range("H131").Value = 10417824
$ws.Range("I131").Value = 2238.3333
$ws.Range("J131").Value = 11905765
$ws.Range("K131").Value = 6714.999899999999
$ws.Range("L131").Value = 35717295
$ws.Range("M131").Value = -1674.999899999999
$ws.Range("N131").Value = -35727375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 29.625
$ws.Range("I2").Value = 26.6
$ws.Range("J2").Value = 34.666668
$ws.Range("K2").Value = 26.6
$ws.Range("L2").Value = 34.666668
$ws.Range("M2").Value = 86.40000000000001
$ws.Range("N2").Value = -260.666668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2090.7742
$ws.Range("I68").Value = 1995.5555
$ws.Range("J68").Value = 2222.6155
$ws.Range("K68").Value = 1995.5555
$ws.Range("L68").Value = 2222.6155
$ws.Range("M68").Value = -1246.5555
$ws.Range("N68").Value = -3720.6155

$ws.Range("H71").Value = 2090.7742
$ws.Range("I71").Value = 1995.5555
$ws.Range("J71").Value = 2222.6155
$ws.Range("K71").Value = 9977.7775
$ws.Range("L71").Value = 11113.0775
$ws.Range("M71").Value = -6233.7775
$ws.Range("N71").Value = -18601.0775

$ws.Range("H100").Value = 3149.1365
$ws.Range("I100").Value = 2526
$ws.Range("J100").Value = 3772.2727
$ws.Range("K100").Value = 2526
$ws.Range("L100").Value = 3772.2727
$ws.Range("M100").Value = -1985
$ws.Range("N100").Value = -4854.2727

$ws.Range("H132").Value = 6282.9375
$ws.Range("I132").Value = 7218.8184
$ws.Range("J132").Value = 4224
$ws.Range("K132").Value = 21656.4552
$ws.Range("L132").Value = 12672
$ws.Range("M132").Value = -19126.4552
$ws.Range("N132").Value = -17732

$ws.Range("H136").Value = 3201.9
$ws.Range("I136").Value = 1880.0358
$ws.Range("J136").Value = 6286.25
$ws.Range("K136").Value = 5640.107400000001
$ws.Range("L136").Value = 18858.75
$ws.Range("M136").Value = -3090.107400000001
$ws.Range("N136").Value = -23958.75

$ws.Range("H140").Value = 62572
$ws.Range("I140").Value = 57980
$ws.Range("J140").Value = 63720
$ws.Range("K140").Value = 57980
$ws.Range("L140").Value = 63720
$ws.Range("M140").Value = -52800
$ws.Range("N140").Value = -74080

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 21431414
$ws.Range("I132").Value = 29412890
$ws.Range("J132").Value = 7448.0527
$ws.Range("K132").Value = 88238670
$ws.Range("L132").Value = 22344.1581
$ws.Range("M132").Value = -88236140
$ws.Range("N132").Value = -27404.1581

$ws.Range("H136").Value = 1280.091
$ws.Range("I136").Value = 756.58826
$ws.Range("J136").Value = 3060
$ws.Range("K136").Value = 2269.76478
$ws.Range("L136").Value = 9180
$ws.Range("M136").Value = 280.23522
$ws.Range("N136").Value = -14280

$ws.Range("H138").Value = 69783.336
$ws.Range("I138").Value = 35000
$ws.Range("J138").Value = 87175
$ws.Range("K138").Value = 35000
$ws.Range("L138").Value = 87175
$ws.Range("N138").Value = -97455
$ws.Range("M138").Value = -29860
